# Auto-generated edit script applying the profit-recalculation diff
# to the Typhon_Profits workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1720.5505   # was 1825.5165
$ws.Range("I15").Value = 1720.5505   # was 1825.5165
$ws.Range("K15").Value = 5161.6515   # was 5476.5495
$ws.Range("M15").Value = -4992.6515   # was -5307.5495
$ws.Range("H64").Value = 4414.2856   # was 4100
$ws.Range("I64").Value = 3500   # was 3250
$ws.Range("J64").Value = 4566.6665   # was 4342.857
$ws.Range("K64").Value = 3500   # was 3250
$ws.Range("L64").Value = 4566.6665   # was 4342.857
$ws.Range("M64").Value = -3252   # was -3002
$ws.Range("N64").Value = -5062.6665   # was -4838.857
$ws.Range("H67").Value = 4414.2856   # was 4100
$ws.Range("I67").Value = 3500   # was 3250
$ws.Range("J67").Value = 4566.6665   # was 4342.857
$ws.Range("K67").Value = 3500   # was 3250
$ws.Range("L67").Value = 4566.6665   # was 4342.857
$ws.Range("M67").Value = -2642   # was -2392
$ws.Range("N67").Value = -6282.6665   # was -6058.857
$ws.Range("H74").Value = 5375   # was 5250
$ws.Range("J74").Value = 5375   # was 5250
$ws.Range("L74").Value = 5375   # was 5250
$ws.Range("N74").Value = -7247   # was -7122
$ws.Range("H76").Value = 3480   # was 3444.4443
$ws.Range("I76").Value = 3477.7778   # was 3428.5715
$ws.Range("J76").Value = 3483.3333   # was 3500
$ws.Range("K76").Value = 3477.7778   # was 3428.5715
$ws.Range("L76").Value = 3483.3333   # was 3500
$ws.Range("M76").Value = -3162.7778   # was -3113.5715
$ws.Range("N76").Value = -4113.3333   # was -4130
$ws.Range("H77").Value = 5375   # was 5250
$ws.Range("J77").Value = 5375   # was 5250
$ws.Range("L77").Value = 26875   # was 26250
$ws.Range("N77").Value = -36235   # was -35610
$ws.Range("H79").Value = 3480   # was 3444.4443
$ws.Range("I79").Value = 3477.7778   # was 3428.5715
$ws.Range("J79").Value = 3483.3333   # was 3500
$ws.Range("K79").Value = 3477.7778   # was 3428.5715
$ws.Range("L79").Value = 3483.3333   # was 3500
$ws.Range("M79").Value = -2385.7778   # was -2336.5715
$ws.Range("N79").Value = -5667.3333   # was -5684
$ws.Range("H111").Value = 6229.25   # was 5583.4
$ws.Range("J111").Value = 0   # was 3000
$ws.Range("L111").Value = 0   # was 9000
$ws.Range("N111").ClearContents()   # was -15134
$ws.Range("H113").Value = 62504576   # was 55560028
$ws.Range("I113").Value = 76926510   # was 90912610
$ws.Range("J113").Value = 9500   # was 5972.2856
$ws.Range("K113").Value = 76926510   # was 90912610
$ws.Range("L113").Value = 9500   # was 5972.2856
$ws.Range("M113").Value = -76923256   # was -90909356
$ws.Range("N113").Value = -16008   # was -12480.2856
$ws.Range("H116").Value = 5533   # was 5299.615
$ws.Range("I116").Value = 3000   # was 2899.8
$ws.Range("K116").Value = 3000   # was 2899.8
$ws.Range("M116").Value = 442   # was 542.1999999999998
$ws.Range("H129").Value = 145839.6   # was 179522.72
$ws.Range("J129").Value = 159702.23   # was 201031.62
$ws.Range("L129").Value = 479106.6900000001   # was 603094.86
$ws.Range("N129").Value = -489106.6900000001   # was -613094.86
$ws.Range("H137").Value = 1352.5807   # was 1379.7142
$ws.Range("I137").Value = 1365.2   # was 1397.0435
$ws.Range("K137").Value = 4095.6   # was 4191.1305
$ws.Range("M137").Value = -1545.6   # was -1641.1305

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2242.8   # was 2494.85
$ws.Range("I32").Value = 1958.3191   # was 2226.4575
$ws.Range("K32").Value = 1958.3191   # was 2226.4575
$ws.Range("M32").Value = -1671.3191   # was -1939.4575
$ws.Range("H63").Value = 2908.6667   # was 2950
$ws.Range("I63").Value = 2919   # was 2950
$ws.Range("J63").Value = 2888   # was 0
$ws.Range("K63").Value = 2919   # was 2950
$ws.Range("L63").Value = 2888   # was 0
$ws.Range("M63").Value = -2233   # was -2264
$ws.Range("N63").Value = -4260   # was None
$ws.Range("H66").Value = 2908.6667   # was 2950
$ws.Range("I66").Value = 2919   # was 2950
$ws.Range("J66").Value = 2888   # was 0
$ws.Range("K66").Value = 14595   # was 14750
$ws.Range("L66").Value = 14440   # was 0
$ws.Range("M66").Value = -11163   # was -11318
$ws.Range("N66").Value = -21304   # was None
$ws.Range("H88").Value = 1000054   # was 500977
$ws.Range("J88").Value = 1000054   # was 500977
$ws.Range("L88").Value = 1000054   # was 500977
$ws.Range("N88").Value = -1000866   # was -501789
$ws.Range("H91").Value = 1000054   # was 500977
$ws.Range("J91").Value = 1000054   # was 500977
$ws.Range("L91").Value = 1000054   # was 500977
$ws.Range("N91").Value = -1002862   # was -503785
$ws.Range("H134").Value = 59000   # was 59500
$ws.Range("J134").Value = 59000   # was 59500
$ws.Range("L134").Value = 59000   # was 59500
$ws.Range("N134").Value = -69140   # was -69640

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1562.4   # was 2118.6365
$ws.Range("I20").Value = 1446.3077   # was 2353.6428
$ws.Range("J20").Value = 1778   # was 1707.375
$ws.Range("K20").Value = 1446.3077   # was 2353.6428
$ws.Range("L20").Value = 1778   # was 1707.375
$ws.Range("M20").Value = -1199.3077   # was -2106.6428
$ws.Range("N20").Value = -2272   # was -2201.375
$ws.Range("H94").Value = 804.7619   # was 753.087
$ws.Range("I94").Value = 629.2308   # was 630.0769
$ws.Range("J94").Value = 1090   # was 913
$ws.Range("K94").Value = 629.2308   # was 630.0769
$ws.Range("L94").Value = 1090   # was 913
$ws.Range("M94").Value = -178.2308   # was -179.0769
$ws.Range("N94").Value = -1992   # was -1815

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2500   # was 1014.3333
$ws.Range("I16").Value = 0   # was 879.2727
$ws.Range("K16").Value = 0   # was 879.2727
$ws.Range("M16").ClearContents()   # was -592.2727
$ws.Range("H58").Value = 28053.525   # was 24563.137
$ws.Range("I58").Value = 1834.3636   # was 1878.909
$ws.Range("J58").Value = 64104.875   # was 47247.363
$ws.Range("K58").Value = 1834.3636   # was 1878.909
$ws.Range("L58").Value = 64104.875   # was 47247.363
$ws.Range("M58").Value = -1631.3636   # was -1675.909
$ws.Range("N58").Value = -64510.875   # was -47653.363
$ws.Range("H62").Value = 52635604   # was 52635740
$ws.Range("J62").Value = 5281   # was 5801.2
$ws.Range("L62").Value = 5281   # was 5801.2
$ws.Range("N62").Value = -6529   # was -7049.2
$ws.Range("H65").Value = 52635604   # was 52635740
$ws.Range("J65").Value = 5281   # was 5801.2
$ws.Range("L65").Value = 26405   # was 29006
$ws.Range("N65").Value = -32645   # was -35246
$ws.Range("H113").Value = 2500   # was 1014.3333
$ws.Range("I113").Value = 0   # was 879.2727
$ws.Range("K113").Value = 0   # was 879.2727
$ws.Range("M113").ClearContents()   # was 1290.7273
$ws.Range("H136").Value = 28053.525   # was 24563.137
$ws.Range("I136").Value = 1834.3636   # was 1878.909
$ws.Range("J136").Value = 64104.875   # was 47247.363
$ws.Range("K136").Value = 5503.0908   # was 5636.727000000001
$ws.Range("L136").Value = 192314.625   # was 141742.089
$ws.Range("M136").Value = -2953.0908   # was -3086.727000000001
$ws.Range("N136").Value = -197414.625   # was -146842.089

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 71428700   # was 62500130
$ws.Range("J38").Value = 125000180   # was 100000160
$ws.Range("L38").Value = 375000540   # was 300000480
$ws.Range("N38").Value = -375001234   # was -300001174
$ws.Range("H131").Value = 750.67   # was 753.6900000000001
$ws.Range("J131").Value = 750.67   # was 753.6900000000001
$ws.Range("L131").Value = 2252.01   # was 2261.07
$ws.Range("N131").Value = -12332.01   # was -12341.07

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3777.75   # was 3250.2
$ws.Range("I80").Value = 3465.8333   # was 2990
$ws.Range("J80").Value = 3964.9   # was 3423.6667
$ws.Range("K80").Value = 3465.8333   # was 2990
$ws.Range("L80").Value = 3964.9   # was 3423.6667
$ws.Range("M80").Value = -2467.8333   # was -1992
$ws.Range("N80").Value = -5960.9   # was -5419.6667
$ws.Range("H83").Value = 3777.75   # was 3250.2
$ws.Range("I83").Value = 3465.8333   # was 2990
$ws.Range("J83").Value = 3964.9   # was 3423.6667
$ws.Range("K83").Value = 17329.1665   # was 14950
$ws.Range("L83").Value = 19824.5   # was 17118.3335
$ws.Range("M83").Value = -12337.1665   # was -9958
$ws.Range("N83").Value = -29808.5   # was -27102.3335

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 292694.25   # was 241158.4
$ws.Range("I2").Value = 437593.88   # was 420000.4
$ws.Range("J2").Value = 2895   # was 62316.4
$ws.Range("K2").Value = 437593.88   # was 420000.4
$ws.Range("L2").Value = 2895   # was 62316.4
$ws.Range("M2").Value = -437481.88   # was -419888.4
$ws.Range("N2").Value = -3119   # was -62540.4
$ws.Range("H68").Value = 0   # was 2766
$ws.Range("J68").Value = 0   # was 2766
$ws.Range("L68").Value = 0   # was 2766
$ws.Range("N68").ClearContents()   # was -4264
$ws.Range("H71").Value = 0   # was 2766
$ws.Range("J71").Value = 0   # was 2766
$ws.Range("L71").Value = 0   # was 13830
$ws.Range("N71").ClearContents()   # was -21318
$ws.Range("H82").Value = 0   # was 747.0909
$ws.Range("I82").Value = 0   # was 721.8
$ws.Range("J82").Value = 0   # was 1000
$ws.Range("K82").Value = 0   # was 721.8
$ws.Range("L82").Value = 0   # was 1000
$ws.Range("M82").ClearContents()   # was -360.8
$ws.Range("N82").ClearContents()   # was -1722
$ws.Range("H85").Value = 0   # was 747.0909
$ws.Range("I85").Value = 0   # was 721.8
$ws.Range("J85").Value = 0   # was 1000
$ws.Range("K85").Value = 0   # was 721.8
$ws.Range("L85").Value = 0   # was 1000
$ws.Range("M85").ClearContents()   # was 526.2
$ws.Range("N85").ClearContents()   # was -3496

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 58081444   # was 55024532
$ws.Range("J107").Value = 5682993   # was 5051566
$ws.Range("L107").Value = 17048979   # was 15154698
$ws.Range("N107").Value = -17052819   # was -15158538
